# Automatic update of files.
# Re-shuffles species observation records across rows 47-51, updating
# the Id (A), Taxonsorteringsordning (B), TaxonId (E), Artnamn (F),
# Vetenskapligt namn (G), Auktor (H), Ost (Q), Nord (R) and, where noted,
# the Biotop-beskrivning (AI) columns so each record lands on its new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 (was id 111974185 / Blå taggsvamp) -> now id 111974186 / Skrovlig taggsvamp
$ws.Cells.Item(47, 1).Value  = 111974186      # A47
$ws.Cells.Item(47, 2).Value  = 90816          # B47
$ws.Cells.Item(47, 5).Value  = 2059           # E47
$ws.Cells.Item(47, 6).Value  = "Skrovlig taggsvamp"                         # F47
$ws.Cells.Item(47, 7).Value  = "Hydnellum scabrosum"                        # G47
$ws.Cells.Item(47, 8).Value  = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"       # H47
$ws.Cells.Item(47, 17).Value = 439860         # Q47
$ws.Cells.Item(47, 18).Value = 6952250        # R47

# Row 48 (was id 111974191 / Talltaggsvamp) -> now id 111974187 / Svart taggsvamp
$ws.Cells.Item(48, 1).Value  = 111974187      # A48
$ws.Cells.Item(48, 2).Value  = 90844          # B48
$ws.Cells.Item(48, 5).Value  = 5449           # E48
$ws.Cells.Item(48, 6).Value  = "Svart taggsvamp"                            # F48
$ws.Cells.Item(48, 7).Value  = "Phellodon niger"                            # G48
$ws.Cells.Item(48, 8).Value  = "(Fr.:Fr.) P.Karst."                         # H48
$ws.Cells.Item(48, 17).Value = 439865         # Q48
$ws.Cells.Item(48, 18).Value = 6952242        # R48
$ws.Cells.Item(48, 35).Value = "äldre renbetad ristallskog med lavfläckar på torr moränmark, under tallåga"  # AI48

# Row 49 (was id 111974186 / Skrovlig taggsvamp) -> now id 111974188 / Talltaggsvamp
$ws.Cells.Item(49, 1).Value  = 111974188      # A49
$ws.Cells.Item(49, 2).Value  = 90786          # B49
$ws.Cells.Item(49, 5).Value  = 3100           # E49
$ws.Cells.Item(49, 6).Value  = "Talltaggsvamp"                              # F49
$ws.Cells.Item(49, 7).Value  = "Bankera fuligineoalba"                      # G49
$ws.Cells.Item(49, 8).Value  = "(Schmidt : Fr.) Pouzar"                     # H49
$ws.Cells.Item(49, 17).Value = 439870         # Q49
$ws.Cells.Item(49, 18).Value = 6952225        # R49

# Row 50 (was id 111974187 / Svart taggsvamp) -> now id 111974185 / Blå taggsvamp
$ws.Cells.Item(50, 1).Value  = 111974185      # A50
$ws.Cells.Item(50, 2).Value  = 90794          # B50
$ws.Cells.Item(50, 5).Value  = 4362           # E50
$ws.Cells.Item(50, 6).Value  = "Blå taggsvamp"                              # F50
$ws.Cells.Item(50, 7).Value  = "Hydnellum caeruleum"                        # G50
$ws.Cells.Item(50, 8).Value  = "(Hornem.) P.Karst."                         # H50
$ws.Cells.Item(50, 17).Value = 439827         # Q50
$ws.Cells.Item(50, 18).Value = 6952233        # R50
$ws.Cells.Item(50, 35).Value = "äldre renbetad ristallskog med lavfläckar på torr moränmark"  # AI50

# Row 51 (was id 111974188 / Talltaggsvamp) -> now id 111974191 / Talltaggsvamp
$ws.Cells.Item(51, 1).Value  = 111974191      # A51
$ws.Cells.Item(51, 2).Value  = 90786          # B51
$ws.Cells.Item(51, 17).Value = 439978         # Q51
$ws.Cells.Item(51, 18).Value = 6952214        # R51
